$d = $word.ActiveDocument

# Map of old exact text -> new exact text, applied as whole-field replacements
# (date header + each of the 25 "NNN÷N=NN, N" table-cell answers).
$replacements = @(
    @("2025-07-20 Sunday", "2025-07-21 Monday"),
    @("173÷7=24, 5", "408÷8=51, 0"),
    @("328÷2=164, 0", "324÷2=162, 0"),
    @("941÷9=104, 5", "634÷7=90, 4"),
    @("880÷6=146, 4", "716÷7=102, 2"),
    @("643÷7=91, 6", "513÷4=128, 1"),
    @("726÷7=103, 5", "181÷4=45, 1"),
    @("227÷5=45, 2", "615÷9=68, 3"),
    @("107÷3=35, 2", "708÷4=177, 0"),
    @("986÷8=123, 2", "480÷7=68, 4"),
    @("930÷5=186, 0", "512÷7=73, 1"),
    @("824÷7=117, 5", "883÷7=126, 1"),
    @("827÷6=137, 5", "503÷4=125, 3"),
    @("753÷4=188, 1", "281÷6=46, 5"),
    @("497÷2=248, 1", "302÷5=60, 2"),
    @("220÷8=27, 4", "473÷6=78, 5"),
    @("269÷4=67, 1", "345÷7=49, 2"),
    @("681÷5=136, 1", "916÷7=130, 6"),
    @("364÷2=182, 0", "376÷7=53, 5"),
    @("804÷3=268, 0", "242÷2=121, 0"),
    @("453÷4=113, 1", "947÷5=189, 2"),
    @("175÷5=35, 0", "732÷7=104, 4"),
    @("437÷9=48, 5", "300÷7=42, 6"),
    @("448÷3=149, 1", "916÷7=130, 6"),
    @("425÷3=141, 2", "524÷7=74, 6"),
    @("509÷9=56, 5", "848÷7=121, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "FAILED to replace: $old -> $new"
    }
}
